$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
$ws2 = $wb.Worksheets.Item("Add Devices Loop B")

# ---------------------------------------------------------------------------
# New label/value text introduced by the "loop loading details" columns.
# ---------------------------------------------------------------------------
$voltDrop = "Volt Drop (V)"
$voltDropWorst = "Volt Drop (worst case)"
$loadingDetailHeader = "Volt drop loading detail name"
$loadingDetailWorstHeader = "Volt drop worst case loading detail name"

foreach ($ws in @($ws1, $ws2)) {

    # Row heights for the (now taller / wrapped) constants rows.
    $ws.Rows.Item(3).RowHeight = 28.8
    $ws.Rows.Item(4).RowHeight = 43.2

    # New "loading detail name" constant labels (E3/E4), styled like the
    # other grey constant boxes (grey fill + border) but left aligned and
    # wrapping, matching the "Volt Drop (V)" / "Volt Drop (worst case)"
    # values used elsewhere.
    $ws.Range("A6").Copy()
    $ws.Range("E3").PasteSpecial(-4122)
    $ws.Range("E3").Value = $voltDrop
    $ws.Range("E3").HorizontalAlignment = -4131
    $ws.Range("E3").WrapText = $true

    $ws.Range("A6").Copy()
    $ws.Range("E4").PasteSpecial(-4122)
    $ws.Range("E4").Value = $voltDropWorst
    $ws.Range("E4").HorizontalAlignment = -4131
    $ws.Range("E4").WrapText = $true

    # New header cells for the two extra columns, styled like the other bold
    # bordered headers in row 5 (A5:C5, F5:G5).
    $ws.Range("F5").Copy()
    $ws.Range("H5").PasteSpecial(-4122)
    $ws.Range("H5").Value = $loadingDetailHeader

    $ws.Range("F5").Copy()
    $ws.Range("I5").PasteSpecial(-4122)
    $ws.Range("I5").Value = $loadingDetailWorstHeader

    # New data cells for row 6, styled like the new E3/E4 constant boxes.
    $ws.Range("E3").Copy()
    $ws.Range("H6").PasteSpecial(-4122)
    $ws.Range("H6").Value = $voltDrop

    $ws.Range("E3").Copy()
    $ws.Range("I6").PasteSpecial(-4122)
    $ws.Range("I6").Value = $voltDropWorst
}

# Sheet1 ("Add Devices Loop A") has an extra data row (row 7) that also needs
# the new columns.
$ws1.Range("E3").Copy()
$ws1.Range("H7").PasteSpecial(-4122)
$ws1.Range("H7").Value = $voltDrop

$ws1.Range("E3").Copy()
$ws1.Range("I7").PasteSpecial(-4122)
$ws1.Range("I7").Value = $voltDropWorst

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selection / active-tab swap: "Add Devices Loop A" becomes the active tab,
# and both sheets now have H5:I6 selected (instead of F6:G6).
# ---------------------------------------------------------------------------
$ws2.Range("H5:I6").Select()
$ws1.Activate()
$ws1.Range("H5:I6").Select()
